{"js": "// Fix HTML original/edited handling: proper UTF-8 en-dash markers and\n// substitution; update DOCX markup en-dash parsing.\n//\n// Four runs in the document hold a \"raw markup\" preview of the\n// Critic-markdown syntax. They need to be normalized so the en-dash\n// wrapper markers (and, for the replacement example, the \"~>replacement\"\n// tail) are stripped, leaving just the plain original text.\n\nconst EN_DASH = \"\\u2013\"; // \u2013\n\nconst replacements = [\n  {\n    find: \" {\" + EN_DASH + \" removed\" + EN_DASH + \"}\",\n    replace: \"  removed\",\n  },\n  {\n    find: \" {\" + EN_DASH + \"This is removed text\" + EN_DASH + \"}\",\n    replace: \" This is removed text\",\n  },\n  {\n    find: \" {\" + EN_DASH + \"This is removed text with some added\" + EN_DASH + \"}\",\n    replace: \" This is removed text with some added\",\n  },\n  {\n    find: \" {This is original text~> this is the replacement}\",\n    replace: \" This is original text\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Fix HTML original/edited handling: proper UTF-8 en-dash markers and\n# substitution; update DOCX markup en-dash parsing.\n#\n# Four runs hold a \"raw markup\" preview of the Critic-markdown syntax.\n# Normalize them so the en-dash wrapper markers (and, for the replacement\n# example, the \"~>replacement\" tail) are stripped, leaving just the plain\n# original text.\n\n$d = $word.ActiveDocument\n\n$enDash = [char]0x2013\n\n$replacements = @(\n    @{ Find = \" {$enDash removed$enDash}\"; Replace = \"  removed\" },\n    @{ Find = \" {$enDash\" + \"This is removed text$enDash}\"; Replace = \" This is removed text\" },\n    @{ Find = \" {$enDash\" + \"This is removed text with some added$enDash}\"; Replace = \" This is removed text with some added\" },\n    @{ Find = \" {This is original text~> this is the replacement}\"; Replace = \" This is original text\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
